# Apply updated Price (D) / Volume(1h) (E) values for cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($cell, [string]$val)
    # Force the string to be stored as text, never auto-coerced to a number/date,
    # by temporarily switching the cell to the Text number format, then restoring
    # its original style so no visible formatting changes remain.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextCellValue $ws.Cells.Item(2, 4) '26.267.07'
Set-TextCellValue $ws.Cells.Item(2, 5) '  +1.03%  '
Set-TextCellValue $ws.Cells.Item(3, 4) '1.625.51'
Set-TextCellValue $ws.Cells.Item(3, 5) '  +1.49%  '
Set-TextCellValue $ws.Cells.Item(4, 5) '  +0.16%  '
Set-TextCellValue $ws.Cells.Item(5, 4) '212.73'
Set-TextCellValue $ws.Cells.Item(5, 5) '  +0.47%  '
Set-TextCellValue $ws.Cells.Item(6, 5) '  +0.14%  '
Set-TextCellValue $ws.Cells.Item(7, 5) '  +0.39%  '
Set-TextCellValue $ws.Cells.Item(8, 5) '  +0.73%  '
Set-TextCellValue $ws.Cells.Item(9, 4) '0.0616'
Set-TextCellValue $ws.Cells.Item(9, 5) '  +0.42%  '
Set-TextCellValue $ws.Cells.Item(10, 4) '18.98'
Set-TextCellValue $ws.Cells.Item(10, 5) '  +4.67%  '
Set-TextCellValue $ws.Cells.Item(11, 4) '0.0819'
Set-TextCellValue $ws.Cells.Item(11, 5) '  +0.95%  '
Set-TextCellValue $ws.Cells.Item(12, 4) '1.851.72'
Set-TextCellValue $ws.Cells.Item(12, 5) '  +1.47%  '
Set-TextCellValue $ws.Cells.Item(13, 4) '1.621.42'
Set-TextCellValue $ws.Cells.Item(13, 5) '  +1.25%  '
Set-TextCellValue $ws.Cells.Item(14, 4) '4.04'
Set-TextCellValue $ws.Cells.Item(14, 5) '  +0.78%  '
Set-TextCellValue $ws.Cells.Item(15, 4) '0.520'
Set-TextCellValue $ws.Cells.Item(15, 5) '  +1.46%  '
Set-TextCellValue $ws.Cells.Item(16, 4) '26.290.29'
Set-TextCellValue $ws.Cells.Item(16, 5) '  +1.08%  '
Set-TextCellValue $ws.Cells.Item(17, 4) '62.61'
Set-TextCellValue $ws.Cells.Item(17, 5) '  +3.91%  '
Set-TextCellValue $ws.Cells.Item(18, 4) '0.0₃0730'
Set-TextCellValue $ws.Cells.Item(18, 5) '  +1.00%  '
Set-TextCellValue $ws.Cells.Item(19, 5) '  +0.13%  '
Set-TextCellValue $ws.Cells.Item(20, 4) '203.87'
Set-TextCellValue $ws.Cells.Item(20, 5) '  +0.57%  '
Set-TextCellValue $ws.Cells.Item(21, 4) '4.30'
Set-TextCellValue $ws.Cells.Item(21, 5) '  +1.41%  '
Set-TextCellValue $ws.Cells.Item(22, 4) '9.37'
Set-TextCellValue $ws.Cells.Item(22, 5) '  +1.03%  '
Set-TextCellValue $ws.Cells.Item(23, 4) '6.05'
Set-TextCellValue $ws.Cells.Item(23, 5) '  +0.72%  '
Set-TextCellValue $ws.Cells.Item(24, 5) '  +7.60%  '
Set-TextCellValue $ws.Cells.Item(25, 4) '143.54'
Set-TextCellValue $ws.Cells.Item(25, 5) '  +1.29%  '
Set-TextCellValue $ws.Cells.Item(26, 5) '  +0.16%  '
Set-TextCellValue $ws.Cells.Item(27, 4) '0.121'
Set-TextCellValue $ws.Cells.Item(27, 5) '  -0.43%  '
Set-TextCellValue $ws.Cells.Item(28, 4) '15.30'
Set-TextCellValue $ws.Cells.Item(28, 5) '  +1.14%  '
Set-TextCellValue $ws.Cells.Item(29, 4) '6.57'
Set-TextCellValue $ws.Cells.Item(29, 5) '  +2.23%  '
Set-TextCellValue $ws.Cells.Item(30, 4) '0.0527'
Set-TextCellValue $ws.Cells.Item(30, 5) '  +10.67%  '
Set-TextCellValue $ws.Cells.Item(31, 5) '  +0.66%  '
Set-TextCellValue $ws.Cells.Item(32, 4) '3.19'
Set-TextCellValue $ws.Cells.Item(32, 5) '  +2.58%  '
Set-TextCellValue $ws.Cells.Item(33, 4) '2.96'
Set-TextCellValue $ws.Cells.Item(33, 5) '  -0.04%  '
Set-TextCellValue $ws.Cells.Item(34, 5) '  +2.32%  '
Set-TextCellValue $ws.Cells.Item(35, 5) '  +1.02%  '
Set-TextCellValue $ws.Cells.Item(36, 4) '1.171.61'
Set-TextCellValue $ws.Cells.Item(36, 5) '  +3.76%  '
Set-TextCellValue $ws.Cells.Item(37, 4) '0.0165'
Set-TextCellValue $ws.Cells.Item(37, 5) '  +2.17%  '
Set-TextCellValue $ws.Cells.Item(38, 4) '0.809'
Set-TextCellValue $ws.Cells.Item(38, 5) '  +2.25%  '
Set-TextCellValue $ws.Cells.Item(39, 5) '  +0.14%  '
Set-TextCellValue $ws.Cells.Item(40, 5) '  +0.57%  '
Set-TextCellValue $ws.Cells.Item(41, 4) '0.499'
Set-TextCellValue $ws.Cells.Item(41, 5) '  +1.64%  '
Set-TextCellValue $ws.Cells.Item(42, 4) '0.793'
Set-TextCellValue $ws.Cells.Item(42, 5) '  +1.05%  '
Set-TextCellValue $ws.Cells.Item(43, 4) '5.30'
Set-TextCellValue $ws.Cells.Item(43, 5) '  +3.08%  '
Set-TextCellValue $ws.Cells.Item(44, 4) '1.762.85'
Set-TextCellValue $ws.Cells.Item(44, 5) '  +1.52%  '
Set-TextCellValue $ws.Cells.Item(45, 4) '93.42'
Set-TextCellValue $ws.Cells.Item(45, 5) '  +0.54%  '
Set-TextCellValue $ws.Cells.Item(46, 5) '  +14.42%  '
Set-TextCellValue $ws.Cells.Item(47, 4) '1.52'
Set-TextCellValue $ws.Cells.Item(47, 5) '  +0.29%  '
Set-TextCellValue $ws.Cells.Item(48, 4) '54.27'
Set-TextCellValue $ws.Cells.Item(48, 5) '  +1.43%  '
Set-TextCellValue $ws.Cells.Item(49, 4) '0.0509'
Set-TextCellValue $ws.Cells.Item(49, 5) '  +1.08%  '
Set-TextCellValue $ws.Cells.Item(50, 4) '0.409'
Set-TextCellValue $ws.Cells.Item(50, 5) '  +0.44%  '
Set-TextCellValue $ws.Cells.Item(51, 5) '  +0.04%  '
